$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 687.1667
$ws.Range("I9").Value = 223.66667
$ws.Range("J9").Value = 1150.6666
$ws.Range("K9").Value = 223.66667
$ws.Range("L9").Value = 1150.6666
$ws.Range("M9").Value = -54.66667000000001
$ws.Range("N9").Value = -1488.6666

$ws.Range("H29").Value = 701
$ws.Range("I29").Value = 103
$ws.Range("K29").Value = 309
$ws.Range("M29").Value = -28

$ws.Range("H43").Value = 29413820
$ws.Range("I43").Value = 71431010
$ws.Range("J43").Value = 1788
$ws.Range("K43").Value = 71431010
$ws.Range("L43").Value = 1788
$ws.Range("M43").Value = -71430941
$ws.Range("N43").Value = -1926

$ws.Range("H62").Value = 22739228
$ws.Range("I62").Value = 35731704
$ws.Range("J62").Value = 2400
$ws.Range("K62").Value = 35731704
$ws.Range("L62").Value = 2400
$ws.Range("M62").Value = -35731080
$ws.Range("N62").Value = -3648

$ws.Range("H65").Value = 22739228
$ws.Range("I65").Value = 35731704
$ws.Range("J65").Value = 2400
$ws.Range("K65").Value = 178658520
$ws.Range("L65").Value = 12000
$ws.Range("M65").Value = -178655400
$ws.Range("N65").Value = -18240

$ws.Range("H132").Value = 7357789
$ws.Range("I132").Value = 7940855
$ws.Range("J132").Value = 11160
$ws.Range("K132").Value = 23822565
$ws.Range("L132").Value = 33480
$ws.Range("M132").Value = -23820035
$ws.Range("N132").Value = -38540

$ws.Range("H135").Value = 1100
$ws.Range("I135").Value = 1183.2354
$ws.Range("J135").Value = 156.66667
$ws.Range("K135").Value = 10649.1186
$ws.Range("L135").Value = 1410.00003
$ws.Range("M135").Value = -8114.1186
$ws.Range("N135").Value = -6480.00003

$ws.Range("H138").Value = 2723.04
$ws.Range("I138").Value = 1048.0286
$ws.Range("J138").Value = 3624.9692
$ws.Range("K138").Value = 3144.0858
$ws.Range("L138").Value = 10874.9076
$ws.Range("M138").Value = 1995.9142
$ws.Range("N138").Value = -21154.9076

$ws.Range("H141").Value = 5408.9443
$ws.Range("I141").Value = 2510.7856
$ws.Range("J141").Value = 15552.5
$ws.Range("K141").Value = 7532.3568
$ws.Range("L141").Value = 46657.5
$ws.Range("M141").Value = -2352.3568
$ws.Range("N141").Value = -57017.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 14723085
$ws.Range("I32").Value = 18194614
$ws.Range("J32").Value = 35845.383
$ws.Range("K32").Value = 18194614
$ws.Range("L32").Value = 35845.383
$ws.Range("M32").Value = -18194327
$ws.Range("N32").Value = -36419.383

$ws.Range("H74").Value = 1745.1897
$ws.Range("I74").Value = 1861.561
$ws.Range("J74").Value = 1464.5294
$ws.Range("K74").Value = 1861.561
$ws.Range("L74").Value = 1464.5294
$ws.Range("M74").Value = -987.5609999999999
$ws.Range("N74").Value = -3212.5294

$ws.Range("H77").Value = 1745.1897
$ws.Range("I77").Value = 1861.561
$ws.Range("J77").Value = 1464.5294
$ws.Range("K77").Value = 9307.805
$ws.Range("L77").Value = 7322.646999999999
$ws.Range("M77").Value = -4939.805
$ws.Range("N77").Value = -16058.647

$ws.Range("H97").Value = 1021.86957
$ws.Range("I97").Value = 885.1429000000001
$ws.Range("J97").Value = 1234.5555
$ws.Range("K97").Value = 885.1429000000001
$ws.Range("L97").Value = 1234.5555
$ws.Range("M97").Value = -389.1429000000001
$ws.Range("N97").Value = -2226.5555

$ws.Range("H122").Value = 10162.714
$ws.Range("I122").Value = 17109.75
$ws.Range("J122").Value = 900
$ws.Range("K122").Value = 51329.25
$ws.Range("L122").Value = 2700
$ws.Range("M122").Value = -48879.25
$ws.Range("N122").Value = -7600

$ws.Range("H132").Value = 14927745
$ws.Range("I132").Value = 18182492
$ws.Range("J132").Value = 10149.75
$ws.Range("K132").Value = 54547476
$ws.Range("L132").Value = 30449.25
$ws.Range("M132").Value = -54544946
$ws.Range("N132").Value = -35509.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 23257778
$ws.Range("I105").Value = 2049.1614
$ws.Range("K105").Value = 2049.1614
$ws.Range("M105").Value = -302.1614

$ws.Range("H134").Value = 1688278
$ws.Range("I134").Value = 3332.577
$ws.Range("J134").Value = 7946646.5
$ws.Range("K134").Value = 9997.731
$ws.Range("L134").Value = 23839939.5
$ws.Range("M134").Value = -7462.731
$ws.Range("N134").Value = -23845009.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 38462630
$ws.Range("I58").Value = 55556400
$ws.Range("J58").Value = 1650.125
$ws.Range("K58").Value = 55556400
$ws.Range("L58").Value = 1650.125
$ws.Range("M58").Value = -55556197
$ws.Range("N58").Value = -2056.125

$ws.Range("H122").Value = 187500800
$ws.Range("I122").Value = 187500800
$ws.Range("K122").Value = 562502400
$ws.Range("M122").Value = -562499950

$ws.Range("H134").Value = 1778.4762
$ws.Range("I134").Value = 1833.2667
$ws.Range("K134").Value = 5499.800099999999
$ws.Range("M134").Value = -2964.800099999999

$ws.Range("H136").Value = 38462630
$ws.Range("I136").Value = 55556400
$ws.Range("J136").Value = 1650.125
$ws.Range("K136").Value = 166669200
$ws.Range("L136").Value = 4950.375
$ws.Range("M136").Value = -166666650
$ws.Range("N136").Value = -10050.375

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H69").Value = 2175.4707
$ws.Range("I69").Value = 900
$ws.Range("J69").Value = 2345.5334
$ws.Range("K69").Value = 2700
$ws.Range("L69").Value = 7036.600199999999
$ws.Range("M69").Value = -1889
$ws.Range("N69").Value = -8658.600199999999

$ws.Range("H72").Value = 2175.4707
$ws.Range("I72").Value = 900
$ws.Range("J72").Value = 2345.5334
$ws.Range("K72").Value = 8100
$ws.Range("L72").Value = 21109.8006
$ws.Range("M72").Value = -4044
$ws.Range("N72").Value = -29221.8006

$ws.Range("H131").Value = 753.13
$ws.Range("J131").Value = 767.69147
$ws.Range("L131").Value = 2303.07441
$ws.Range("N131").Value = -12383.07441

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H7").Value = 3750000.5
$ws.Range("J7").Value = 5000000
$ws.Range("L7").Value = 5000000
$ws.Range("N7").Value = -5000224

$ws.Range("H8").Value = 3750000.5
$ws.Range("J8").Value = 5000000
$ws.Range("L8").Value = 5000000
$ws.Range("N8").Value = -5000278

$ws.Range("H102").Value = 1266.9
$ws.Range("I102").Value = 1137.375
$ws.Range("K102").Value = 1137.375
$ws.Range("M102").Value = 484.625

$ws.Range("H113").Value = 2064.125
$ws.Range("I113").Value = 2100
$ws.Range("J113").Value = 2004.3334
$ws.Range("K113").Value = 2100
$ws.Range("L113").Value = 2004.3334
$ws.Range("M113").Value = 70
$ws.Range("N113").Value = -6344.3334

$ws.Range("H122").Value = 17864436
$ws.Range("I122").Value = 27786500
$ws.Range("J122").Value = 4720.7
$ws.Range("K122").Value = 83359500
$ws.Range("L122").Value = 14162.1
$ws.Range("M122").Value = -83357050
$ws.Range("N122").Value = -19062.1

$ws.Range("H123").Value = 17058.223
$ws.Range("J123").Value = 17058.223
$ws.Range("L123").Value = 17058.223
$ws.Range("N123").Value = -21958.223

$ws.Range("H132").Value = 6683.069
$ws.Range("I132").Value = 1910.5883
$ws.Range("K132").Value = 5731.7649
$ws.Range("M132").Value = -3201.7649

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 31253200
$ws.Range("I40").Value = 3433.3333
$ws.Range("K40").Value = 3433.3333
$ws.Range("M40").Value = -3297.3333

$ws.Range("H61").Value = 1872.0526
$ws.Range("I61").Value = 1562.8572
$ws.Range("K61").Value = 1562.8572
$ws.Range("M61").Value = -1360.8572

$ws.Range("H82").Value = 1616.6666
$ws.Range("I82").Value = 1466.6666
$ws.Range("J82").Value = 1766.6666
$ws.Range("K82").Value = 1466.6666
$ws.Range("L82").Value = 1766.6666
$ws.Range("M82").Value = -1105.6666
$ws.Range("N82").Value = -2488.6666

$ws.Range("H85").Value = 1616.6666
$ws.Range("I85").Value = 1466.6666
$ws.Range("J85").Value = 1766.6666
$ws.Range("K85").Value = 1466.6666
$ws.Range("L85").Value = 1766.6666
$ws.Range("M85").Value = -218.6666
$ws.Range("N85").Value = -4262.6666

$ws.Range("H113").Value = 1872.0526
$ws.Range("I113").Value = 1562.8572
$ws.Range("K113").Value = 1562.8572
$ws.Range("M113").Value = 607.1428000000001

$ws.Range("H122").Value = 50600
$ws.Range("I122").Value = 50600
$ws.Range("K122").Value = 151800
$ws.Range("M122").Value = -149350

$ws.Range("H132").Value = 17862400
$ws.Range("I132").Value = 34484924
$ws.Range("J132").Value = 8581.593000000001
$ws.Range("K132").Value = 103454772
$ws.Range("L132").Value = 25744.779
$ws.Range("M132").Value = -103452242
$ws.Range("N132").Value = -30804.779

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 25599.889
$ws.Range("I70").Value = 12000
$ws.Range("J70").Value = 27299.875
$ws.Range("K70").Value = 12000
$ws.Range("L70").Value = 27299.875
$ws.Range("M70").Value = -11685
$ws.Range("N70").Value = -27929.875

$ws.Range("H73").Value = 25599.889
$ws.Range("I73").Value = 12000
$ws.Range("J73").Value = 27299.875
$ws.Range("K73").Value = 12000
$ws.Range("L73").Value = 27299.875
$ws.Range("M73").Value = -10908
$ws.Range("N73").Value = -29483.875

$ws.Range("H75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").ClearContents()

$ws.Range("H78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").ClearContents()

$ws.Range("H113").Value = 41667130
$ws.Range("I113").Value = 55555916
$ws.Range("J113").Value = 762.1667
$ws.Range("K113").Value = 166667748
$ws.Range("L113").Value = 2286.5001
$ws.Range("M113").Value = -166665578
$ws.Range("N113").Value = -6626.5001

$ws.Range("H122").Value = 3233.6
$ws.Range("I122").Value = 2172
$ws.Range("K122").Value = 6516
$ws.Range("M122").Value = -4066

$ws.Range("H132").Value = 11249830
$ws.Range("I132").Value = 21762530
$ws.Range("J132").Value = 3686.186
$ws.Range("K132").Value = 65287590
$ws.Range("L132").Value = 11058.558
$ws.Range("M132").Value = -65285060
$ws.Range("N132").Value = -16118.558
